$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a redundant "ID" column A (duplicate of column B) with no
# header of its own. Remove it entirely so every following column shifts
# one slot to the left (B->A, C->B, ..., H->G), matching the corrected
# layout used for the AAPL yearly EPS pie chart.
$ws.Columns.Item(1).Delete()
